$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Sales_Territory / Sales_Area mapping blocks ---
# Rows 4-6 currently map to Sales_Territory (ST_*); they should map to Sales_Area (SA_*)
$ws.Range("E4").Value = "Sales_Area"
$ws.Range("F4").Value = "SA_ID"

$ws.Range("E5").Value = "Sales_Area"
$ws.Range("F5").Value = "SA_Name"

$ws.Range("E6").Value = "Sales_Area"
# F6 (MGR_ID) is unchanged

# Rows 8-10 currently map to Sales_Area (SA_*); they should map to Sales_Territory (ST_*)
$ws.Range("E8").Value = "Sales_Territory"
$ws.Range("F8").Value = "ST_ID"

$ws.Range("E9").Value = "Sales_Territory"
$ws.Range("F9").Value = "ST_Name"

$ws.Range("E10").Value = "Sales_Territory"
# F10 (MGR_ID) is unchanged

# --- Add new row 19 (SSDT drawings marker cell) ---
$ws.Range("E19").Value = " "
$ws.Range("E19").Font.Name = "Arial"
$ws.Range("E19").Font.Size = 10
$ws.Range("E19").Font.Bold = $false

# Move the active selection the same way Excel would after editing near the bottom of the sheet
$ws.Range("D20").Select()
